$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.242.43"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.624.77"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "605.76"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "196.31"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.86%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.651"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "54.12"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "4.199.70"
$ws.Range("E14").Value = "  +3.55%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "13.26"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +5.34%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "594.27"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "70.473.23"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "3.621.39"
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  +0.61%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "17.71"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.15"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.93%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "102.51"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  -0.43%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.82"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.63"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "34.10"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.74"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +4.80%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.24"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "12.36"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0910"
$ws.Range("E34").Value = "  +12.26%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "63.26"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "3.924.74"
$ws.Range("E36").Value = "  +5.30%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.17"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.21%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "528.27"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("E39").Value = "  +0.12%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "37.46"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  +0.52%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.56"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("E44").Value = "  -0.24%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.88"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +2.23%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.38"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("E47").Value = "  +0.52%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.65"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +4.73%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.33"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
